$d = $word.ActiveDocument

$replacements = @(
    @{old = "survival"; new = "Survival"},
    @{old = "pclass";   new = "Pclass"},
    @{old = "name";     new = "Name"},
    @{old = "sex";      new = "Sex"},
    @{old = "age";      new = "Age"},
    @{old = "sibsp";    new = "SibSp"},
    @{old = "parch";    new = "Parch"},
    @{old = "ticket";   new = "Ticket"},
    @{old = "fare";     new = "Fare"},
    @{old = "cabin";    new = "Cabin"},
    @{old = "embarked"; new = "Embarked"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Font.Bold = $true
    $rng.Find.Replacement.Font.Bold = $true
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $true, $r.new, 2) | Out-Null
}
